# Update latest output (run 85)
$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Cells.Item(2, 2).Value = 46040.875          # B2
$schedule.Cells.Item(2, 3).Value = 14.5                # C2
$schedule.Cells.Item(2, 4).Value = 54.81               # D2
$schedule.Cells.Item(2, 5).Value = 206.949366          # E2
$schedule.Cells.Item(2, 6).Value = 3.775759277504105   # F2

$schedule.Cells.Item(3, 1).Value = 46040.9375          # A3
$schedule.Cells.Item(3, 3).Value = 4.5                 # C3
$schedule.Cells.Item(3, 4).Value = 17.01               # D3
$schedule.Cells.Item(3, 5).Value = 426.3550785         # E3
$schedule.Cells.Item(3, 6).Value = 25.0649664021164    # F3

$schedule.Cells.Item(4, 5).Value = -34.45792349999999  # E4
$schedule.Cells.Item(4, 6).Value = -1.012872530864197  # F4

# --- Sheet "Detailed" ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Cells.Item(35, 2).Value = -6.60147
$detailed.Cells.Item(36, 2).Value = -6
$detailed.Cells.Item(37, 2).Value = -6
$detailed.Cells.Item(37, 3).Value = "historical"
$detailed.Cells.Item(38, 2).Value = 17.71906
$detailed.Cells.Item(38, 3).Value = "historical"
$detailed.Cells.Item(39, 2).Value = 30.58443
$detailed.Cells.Item(40, 2).Value = 55.38427
$detailed.Cells.Item(41, 2).Value = 57.3
$detailed.Cells.Item(42, 2).Value = 56.52615
$detailed.Cells.Item(42, 5).Value = "ON"
$detailed.Cells.Item(43, 2).Value = 56.98
$detailed.Cells.Item(43, 5).Value = "ON"
$detailed.Cells.Item(45, 5).Value = "OFF"
$detailed.Cells.Item(46, 2).Value = 57.06007
$detailed.Cells.Item(46, 5).Value = "OFF"
$detailed.Cells.Item(47, 2).Value = 49.04675
$detailed.Cells.Item(48, 2).Value = 36.2
$detailed.Cells.Item(49, 2).Value = 36.2
$detailed.Cells.Item(50, 2).Value = 47.31324
$detailed.Cells.Item(51, 2).Value = 56.36108
$detailed.Cells.Item(54, 2).Value = 48.4027
$detailed.Cells.Item(55, 2).Value = 49.72346
$detailed.Cells.Item(56, 2).Value = 50.63512
$detailed.Cells.Item(57, 2).Value = 56.98
$detailed.Cells.Item(59, 2).Value = 59.11391
$detailed.Cells.Item(60, 2).Value = 57.82778
$detailed.Cells.Item(61, 2).Value = 60.21199
$detailed.Cells.Item(65, 2).Value = 4.45932
$detailed.Cells.Item(66, 2).Value = 0.51
$detailed.Cells.Item(67, 2).Value = 0.51
$detailed.Cells.Item(68, 2).Value = -2.83936
$detailed.Cells.Item(69, 2).Value = -6.13408
$detailed.Cells.Item(70, 2).Value = -6.33411
$detailed.Cells.Item(71, 2).Value = -6.12189
$detailed.Cells.Item(72, 2).Value = -6.50816
$detailed.Cells.Item(73, 2).Value = -6.25482
$detailed.Cells.Item(74, 2).Value = -6.07918
$detailed.Cells.Item(75, 2).Value = -6.73139
$detailed.Cells.Item(76, 2).Value = -6.49292
$detailed.Cells.Item(77, 2).Value = -6.21235
$detailed.Cells.Item(78, 2).Value = -5.51
$detailed.Cells.Item(79, 2).Value = -5.51
$detailed.Cells.Item(80, 2).Value = -5.01
$detailed.Cells.Item(81, 2).Value = -0.96252
$detailed.Cells.Item(82, 2).Value = 0.00025
$detailed.Cells.Item(83, 2).Value = -2.48837
$detailed.Cells.Item(84, 2).Value = -0.50876
$detailed.Cells.Item(85, 2).Value = 0.91882
$detailed.Cells.Item(86, 2).Value = 20.22705
$detailed.Cells.Item(87, 2).Value = 53.6751
$detailed.Cells.Item(88, 2).Value = 57.3908
$detailed.Cells.Item(89, 2).Value = 77.49445
$detailed.Cells.Item(90, 2).Value = 75.93411
$detailed.Cells.Item(91, 2).Value = 65
$detailed.Cells.Item(92, 2).Value = 68.44318
$detailed.Cells.Item(94, 2).Value = 62.9353
$detailed.Cells.Item(95, 2).Value = 59.00814
$detailed.Cells.Item(96, 2).Value = 58.18858
$detailed.Cells.Item(97, 2).Value = 61.96632
